$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that sits after "yang ".
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2) Collapse the split-up field placeholder text into one run with the
#    new combined expression.
# ---------------------------------------------------------------------
$lq = [char]0x2018
$rq = [char]0x2019
$newField = "[wbp.no_srt_pmt] [wbp.nmr_srt_thn; when wbp.no_srt_pmt = " + $lq + $rq + "]"
$rng = $d.Content
$rng.Find.Execute("[wbp.no_srt_pmt; ifempty=wbp.nmr_srt_thn]", $true, $false, $false, $false, $false, `
    $true, 1, $false, $newField, 2)

# ---------------------------------------------------------------------
# 3) Resize the 2nd/3rd data columns of the main table.
# ---------------------------------------------------------------------
$t = $d.Tables.Item(2)
$t.Columns.Item(2).Width = 170.1
$t.Columns.Item(3).Width = 198.45
$t.PreferredWidth = 823.55

# ---------------------------------------------------------------------
# 4) Re-insert the "_GoBack" bookmark in its new spot, splitting the
#    "...ahanan tersebut akan..." run right after "tersebu".
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("ahanan tersebu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $rng2.End
$r2 = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $r2)

# ---------------------------------------------------------------------
# 5) Tighten the line spacing of the blank paragraph right after
#    "...mendapat perhatian sebagaimana mestinya." from 360 -> 240.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("mendapat perhatian sebagaimana mestinya.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$p = $rng3.Paragraphs(1)
$nextPara = $p.Next()
$nextPara.Range.ParagraphFormat.LineSpacing = 12
